$d = $word.ActiveDocument

# The document currently ends with a couple of empty paragraphs followed by
# a paragraph containing "Route 53:". That heading (and the two blank
# paragraphs that led up to it) were just scaffolding from an in-progress
# edit and need to be removed, leaving the document ending on the single
# blank paragraph that existed before them.

$paras = $d.Paragraphs
$count = $paras.Count

$lastPara = $paras.Item($count)            # "Route 53:" paragraph
$startPara = $paras.Item($count - 2)        # first of the two blank paragraphs to drop

$deleteRange = $d.Range($startPara.Range.Start, $lastPara.Range.End)
$deleteRange.Delete()
